# Updates cryptos list prices/volumes (and restores correct Uniswap/ShibaInu
# row ordering) to match the latest scrape, per the GitHub Actions commit.
#
# For numeric-looking price strings we briefly force a Text number format
# before assigning the value (then reset the style back to Normal) so Excel
# keeps them as text like "407.49" / "0.0000167" instead of silently
# re-interpreting them as real numbers.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "62.485.83"
$ws.Range("E2").Value = "  +2.19%  "
$ws.Range("D3").Value = "3.433.35"
$ws.Range("E3").Value = "  +3.30%  "
$ws.Range("E4").Value = "  +0.01%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "407.49"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +2.60%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "131.08"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +5.61%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.602"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +3.08%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "1.00"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +0.01%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.696"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +7.18%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.143"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +22.82%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "42.32"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +4.89%  "
$ws.Range("E12").Value = "  +0.62%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "8.52"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +4.83%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "19.93"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +4.76%  "
$ws.Range("D15").Value = "3.430.98"
$ws.Range("E15").Value = "  +0.66%  "
$ws.Range("D16").Value = "62.633.00"
$ws.Range("E16").Value = "  +2.52%  "
$ws.Range("B17").Value = "Uniswap"
$ws.Range("C17").Value = "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "11.56"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +4.84%  "
$ws.Range("B18").Value = "ShibaInu"
$ws.Range("C18").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.0000167"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +33.21%  "
$ws.Range("E19").Value = "  +3.25%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "3.19"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +1.12%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "84.67"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +7.36%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "315.00"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +6.60%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "12.87"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +2.62%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "3.18"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +3.67%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "4.75"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +0.52%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "29.90"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +4.56%  "
$ws.Range("E27").Value = "  +1.21%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "7.79"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +6.35%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.78"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +12.82%  "
$ws.Range("E30").Value = "  +2.25%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "44.15"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +9.75%  "
$ws.Range("E32").Value = "  +2.71%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "11.47"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +2.92%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.999"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +0.03%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.0488"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +3.98%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "51.62"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -0.62%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.997"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -0.18%  "
$ws.Range("E38").Value = "  +4.34%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.320"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +16.68%  "
$ws.Range("E40").Value = "  -0.96%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "143.52"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +5.47%  "
$ws.Range("E42").Value = "  +4.34%  "
$ws.Range("E43").Value = "  +2.84%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "16.99"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +3.57%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "3.92"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +3.49%  "
$ws.Range("E46").Value = "  +0.49%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "21.44"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +2.49%  "
$ws.Range("D48").Value = "2.112.10"
$ws.Range("E48").Value = "  +0.59%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "2.02"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +11.40%  "
$ws.Range("E50").Value = "  -0.35%  "
$ws.Range("E51").Value = "  +33.56%  "
